# Update "想去人数" (number of interested attendees) values in column F
# across the four worksheets, per the source diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1912
$ws.Range("F3").Value  = 30
$ws.Range("F5").Value  = 815
$ws.Range("F7").Value  = 1024
$ws.Range("F8").Value  = 1675
$ws.Range("F9").Value  = 1318
$ws.Range("F10").Value = 1613
$ws.Range("F11").Value = 1626
$ws.Range("F12").Value = 373
$ws.Range("F13").Value = 1742
$ws.Range("F15").Value = 1190
$ws.Range("F16").Value = 60
$ws.Range("F17").Value = 120
$ws.Range("F18").Value = 2141
$ws.Range("F19").Value = 294
$ws.Range("F20").Value = 841
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 1360
$ws.Range("F25").Value = 1139
$ws.Range("F26").Value = 113
$ws.Range("F28").Value = 1255
$ws.Range("F29").Value = 930
$ws.Range("F30").Value = 27
$ws.Range("F31").Value = 1233
$ws.Range("F32").Value = 89
$ws.Range("F33").Value = 1183
$ws.Range("F34").Value = 355
$ws.Range("F35").Value = 105
$ws.Range("F38").Value = 1751
$ws.Range("F39").Value = 400
$ws.Range("F40").Value = 21
$ws.Range("F42").Value = 2130
$ws.Range("F45").Value = 1356
$ws.Range("F47").Value = 821
$ws.Range("F49").Value = 31

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 285
$ws.Range("F12").Value = 47
$ws.Range("F21").Value = 34
$ws.Range("F22").Value = 8
$ws.Range("F27").Value = 267
$ws.Range("F29").Value = 251
$ws.Range("F31").Value = 66
$ws.Range("F33").Value = 35
$ws.Range("F34").Value = 35
$ws.Range("F36").Value = 25
$ws.Range("F44").Value = 71

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 280
$ws.Range("F5").Value  = 2973
$ws.Range("F6").Value  = 4766
$ws.Range("F10").Value = 840
$ws.Range("F11").Value = 504
$ws.Range("F12").Value = 503
$ws.Range("F13").Value = 1227
$ws.Range("F14").Value = 347
$ws.Range("F15").Value = 910

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1912
$ws.Range("F3").Value  = 280
$ws.Range("F4").Value  = 30
$ws.Range("F5").Value  = 4766
$ws.Range("F7").Value  = 840
$ws.Range("F8").Value  = 504
$ws.Range("F10").Value = 503
$ws.Range("F11").Value = 1227
$ws.Range("F12").Value = 1024
$ws.Range("F13").Value = 1675
$ws.Range("F14").Value = 1318
$ws.Range("F15").Value = 1613
$ws.Range("F16").Value = 1626
$ws.Range("F17").Value = 285
$ws.Range("F19").Value = 1742
$ws.Range("F20").Value = 1190
$ws.Range("F22").Value = 910
$ws.Range("F23").Value = 910
$ws.Range("F24").Value = 2141
$ws.Range("F26").Value = 294
$ws.Range("F27").Value = 841
$ws.Range("F30").Value = 1360
$ws.Range("F32").Value = 1139
$ws.Range("F33").Value = 113
$ws.Range("F34").Value = 1255
$ws.Range("F35").Value = 930
$ws.Range("F36").Value = 1233
$ws.Range("F37").Value = 89
$ws.Range("F39").Value = 1183
$ws.Range("F40").Value = 355
$ws.Range("F43").Value = 1751
$ws.Range("F44").Value = 35
$ws.Range("F45").Value = 21
$ws.Range("F46").Value = 2130
$ws.Range("F49").Value = 1357
$ws.Range("F50").Value = 821
$ws.Range("F52").Value = 71
